$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 <- data from original row 16
$ws.Range("A9").Value = 111670575
$ws.Range("B9").Value = 96346
$ws.Range("E9").Value = 620
$ws.Range("Q9").Value = 558082.6649719321
$ws.Range("R9").Value = 7067974.943554637
$ws.Range("D9").Value = "NT"
$ws.Range("F9").Value = "Skogsfru"
$ws.Range("G9").Value = "Epipogium aphyllum"
$ws.Range("H9").Value = "Sw."

# Row 11 <- data from original row 17
$ws.Range("A11").Value = 111671406
$ws.Range("B11").Value = 78578
$ws.Range("E11").Value = 6458
$ws.Range("Q11").Value = 557823.3030943703
$ws.Range("R11").Value = 7068159.357501161
$ws.Range("D11").Value = "NT"
$ws.Range("F11").Value = "Lunglav"
$ws.Range("G11").Value = "Lobaria pulmonaria"
$ws.Range("H11").Value = "(L.) Hoffm."

# Row 12 <- data from original row 18
$ws.Range("A12").Value = 111671345
$ws.Range("B12").Value = 96348
$ws.Range("E12").Value = 220787
$ws.Range("Q12").Value = 557812.5300353739
$ws.Range("R12").Value = 7068166.248475613
$ws.Range("D12").Value = "VU"
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."

# Row 13 <- data from original row 12
$ws.Range("A13").Value = 111671395
$ws.Range("B13").Value = 96348
$ws.Range("E13").Value = 220787
$ws.Range("Q13").Value = 557763.2623863788
$ws.Range("R13").Value = 7068264.582601988
$ws.Range("D13").Value = "VU"
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."

# Row 14 <- data from original row 9
$ws.Range("A14").Value = 111670593
$ws.Range("B14").Value = 78578
$ws.Range("E14").Value = 6458
$ws.Range("Q14").Value = 558040.5475534229
$ws.Range("R14").Value = 7067901.063021242
$ws.Range("D14").Value = "NT"
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."

# Row 15 <- data from original row 14
$ws.Range("A15").Value = 111670607
$ws.Range("B15").Value = 96368
$ws.Range("E15").Value = 221952
$ws.Range("Q15").Value = 558031.5471372061
$ws.Range("R15").Value = 7067907.98648507
$ws.Range("D15").Value = "LC"
$ws.Range("F15").Value = "Spindelblomster"
$ws.Range("G15").Value = "Neottia cordata"
$ws.Range("H15").Value = "(L.) Rich."

# Row 16 <- data from original row 15
$ws.Range("A16").Value = 111670588
$ws.Range("B16").Value = 96348
$ws.Range("E16").Value = 220787
$ws.Range("Q16").Value = 558039.6361001397
$ws.Range("R16").Value = 7067902.375451046
$ws.Range("D16").Value = "VU"
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."

# Row 17 <- data from original row 13
$ws.Range("A17").Value = 111671384
$ws.Range("B17").Value = 96348
$ws.Range("E17").Value = 220787
$ws.Range("Q17").Value = 557798.0632258818
$ws.Range("R17").Value = 7068181.046264404
$ws.Range("D17").Value = "VU"
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."

# Row 18 <- data from original row 11
$ws.Range("A18").Value = 111670599
$ws.Range("B18").Value = 96348
$ws.Range("E18").Value = 220787
$ws.Range("Q18").Value = 558031.5226908802
$ws.Range("R18").Value = 7067909.315233406
$ws.Range("D18").Value = "VU"
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
